$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 9 new data rows (rows 22-30) to the master-reg_center_machine table,
# continuing the regcntr_id / machine_id sequences and repeating the same
# lang_code / is_active / cr_by / cr_dtimes values used throughout the sheet.
$newRows = @(
    @(10002, 10021),
    @(10003, 10022),
    @(10004, 10023),
    @(10005, 10024),
    @(10006, 10025),
    @(10007, 10026),
    @(10008, 10027),
    @(10009, 10028),
    @(10010, 10029)
)

$row = 22
foreach ($pair in $newRows) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $ws.Cells.Item($row, 3).Value = "eng"
    $ws.Cells.Item($row, 4).Value = $true
    $ws.Cells.Item($row, 5).Value = "superadmin"
    $ws.Cells.Item($row, 6).Value = "now()"
    $row++
}

# Match the author's final view state: scrolled down with the row right
# after the new data selected.
$ws.Range("A31:XFD1048576").Select()

# Set the sheet up for portrait printing, as in the saved workbook.
$ws.PageSetup.Orientation = 1
